$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose updated value is numeric-looking as Text so Excel
# stores the literal string (preserving trailing zeros, etc.) instead of
# silently converting it to a Number. Looping (rather than a single comma-
# unioned Range) so the format reliably lands on every one of the cells.
$textFormatCells = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13",
    "D14", "D15", "D16", "D17", "D21", "D22", "D23", "D24", "D25", "D27",
    "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37",
    "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D48", "D49",
    "D50"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values row by row.
# Row 2
$ws.Range("D2").Value = "27.238.95"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("D3").Value = "1.904.72"
$ws.Range("E3").Value = "  +0.14%  "

# Row 4
$ws.Range("D4").Value = "1.002"

# Row 5
$ws.Range("D5").Value = "307.18"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.29%  "

# Row 7
$ws.Range("D7").Value = "0.5360"
$ws.Range("E7").Value = "  +2.76%  "

# Row 8
$ws.Range("D8").Value = "0.3815"
$ws.Range("E8").Value = "  +0.97%  "

# Row 9
$ws.Range("D9").Value = "0.07290"
$ws.Range("E9").Value = "  +0.07%  "

# Row 10
$ws.Range("D10").Value = "22.04"
$ws.Range("E10").Value = "  +3.69%  "

# Row 11
$ws.Range("D11").Value = "0.9026"
$ws.Range("E11").Value = "  -0.16%  "

# Row 12
$ws.Range("D12").Value = "0.08194"
$ws.Range("E12").Value = "  -1.16%  "

# Row 13
$ws.Range("D13").Value = "95.67"
$ws.Range("E13").Value = "  -1.15%  "

# Row 14
$ws.Range("D14").Value = "5.346"
$ws.Range("E14").Value = "  +0.97%  "

# Row 15
$ws.Range("D15").Value = "0.9996"
$ws.Range("E15").Value = "  -0.07%  "

# Row 16
$ws.Range("D16").Value = "14.85"
$ws.Range("E16").Value = "  +1.86%  "

# Row 17
$ws.Range("D17").Value = "0.000008647"
$ws.Range("E17").Value = "  +0.15%  "

# Row 18
$ws.Range("E18").Value = "  +0.20%  "

# Row 19
$ws.Range("D19").Value = "27.279.71"
$ws.Range("E19").Value = "  -0.07%  "

# Row 20
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "1.126.72"
$ws.Range("E20").Value = "  -40.75%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.017"
$ws.Range("E21").Value = "  -1.54%  "

# Row 22
$ws.Range("D22").Value = "10.77"
$ws.Range("E22").Value = "  +0.85%  "

# Row 23
$ws.Range("D23").Value = "6.517"
$ws.Range("E23").Value = "  +1.24%  "

# Row 24
$ws.Range("D24").Value = "149.76"
$ws.Range("E24").Value = "  +1.71%  "

# Row 25
$ws.Range("D25").Value = "2.289"
$ws.Range("E25").Value = "  -1.21%  "

# Row 26
$ws.Range("E26").Value = "  +0.27%  "

# Row 27
$ws.Range("D27").Value = "1.746"
$ws.Range("E27").Value = "  -0.29%  "

# Row 28
$ws.Range("D28").Value = "117.06"
$ws.Range("E28").Value = "  +1.41%  "

# Row 29
$ws.Range("D29").Value = "4.825"
$ws.Range("E29").Value = "  -0.55%  "

# Row 30
$ws.Range("D30").Value = "4.801"
$ws.Range("E30").Value = "  -2.29%  "

# Row 31
$ws.Range("D31").Value = "0.09283"
$ws.Range("E31").Value = "  +0.27%  "

# Row 32
$ws.Range("D32").Value = "0.8356"
$ws.Range("E32").Value = "  +4.23%  "

# Row 33
$ws.Range("D33").Value = "0.05062"
$ws.Range("E33").Value = "  -0.28%  "

# Row 34
$ws.Range("D34").Value = "1.223"
$ws.Range("E34").Value = "  -1.38%  "

# Row 35
$ws.Range("D35").Value = "3.004"
$ws.Range("E35").Value = "  +1.59%  "

# Row 36
$ws.Range("D36").Value = "3.350"
$ws.Range("E36").Value = "  -2.33%  "

# Row 37
$ws.Range("D37").Value = "2.685"
$ws.Range("E37").Value = "  +3.37%  "

# Row 38
$ws.Range("D38").Value = "0.5753"
$ws.Range("E38").Value = "  +0.29%  "

# Row 39
$ws.Range("D39").Value = "0.02005"
$ws.Range("E39").Value = "  -0.04%  "

# Row 40
$ws.Range("D40").Value = "1.077"
$ws.Range("E40").Value = "  -0.29%  "

# Row 41
$ws.Range("D41").Value = "9.338"
$ws.Range("E41").Value = "  +3.38%  "

# Row 42
$ws.Range("D42").Value = "6.563"
$ws.Range("E42").Value = "  -0.42%  "

# Row 43
$ws.Range("D43").Value = "117.68"
$ws.Range("E43").Value = "  +1.21%  "

# Row 44
$ws.Range("D44").Value = "0.1523"
$ws.Range("E44").Value = "  +0.28%  "

# Row 45
$ws.Range("D45").Value = "0.4929"
$ws.Range("E45").Value = "  +1.00%  "

# Row 46
$ws.Range("E46").Value = "  +0.26%  "

# Row 47
$ws.Range("E47").Value = "  +0.66%  "

# Row 48
$ws.Range("D48").Value = "1.640"
$ws.Range("E48").Value = "  +0.64%  "

# Row 49
$ws.Range("D49").Value = "38.56"
$ws.Range("E49").Value = "  +1.56%  "

# Row 50
$ws.Range("D50").Value = "0.06131"
$ws.Range("E50").Value = "  +3.05%  "

# Row 51
$ws.Range("E51").Value = "  -0.99%  "
